$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Decide on project idea" status changes from In Progress -> Complete
$ws.Range("C4").Value = "Complete"

# Append new row 6 (will become row 8 after the later insert) at the bottom first,
# so that the shared string "Complete PA1 powerpoint report" gets registered
# before "Design module outlines" / "Define input data format".
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)
$ws.Range("A6").Value = 45753
$ws.Range("B6").Value = "Complete PA1 powerpoint report"
$ws.Range("C6").Value = "Complete"
$ws.Range("D6").Value = "Everyone"

# Now insert two new blank rows before row 5 ("Design module communication diagram"),
# pushing it (and the row just added) down.
$ws.Range("A5:A6").EntireRow.Insert()

$ws.Range("A5").Value = 45753
$ws.Range("B5").Value = "Design module outlines"
$ws.Range("C5").Value = "Complete"
$ws.Range("D5").Value = "Everyone"

$ws.Range("A6").Value = 45753
$ws.Range("B6").Value = "Define input data format"
$ws.Range("C6").Value = "Complete"
$ws.Range("D6").Value = "Everyone"

# Existing "Design module communication diagram" row is now row 7; update its status
$ws.Range("C7").Value = "Complete"

$ws.Range("B8").Select()
